# Fruta / hortaliza, semanal
#
# A new weekly price report row is inserted at row 313 (pushing the
# existing rows 313..357 down to 314..358), for
# "Macroferia Regional de Talca" / "Betarraga".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 313; this shifts all
# rows from 313 down through 357 to 314 through 358 (the sheet grows
# from A1:R357 to A1:R358).
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with this week's record.
$ws.Range("A313").Value = 5
$ws.Range("B313").Value = "Macroferia Regional de Talca"
$ws.Range("C313").Value = "Maule"
$ws.Range("D313").Value = 44776
$ws.Range("E313").Value = 7
$ws.Range("F313").Value = 100114014
$ws.Range("G313").Value = "Betarraga"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 4000
$ws.Range("K313").Value = 750
$ws.Range("L313").Value = 750
$ws.Range("M313").Value = 750
$ws.Range("N313").Value = "$/paquete 5 unidades"
$ws.Range("O313").Value = "Región del Maule"
$ws.Range("P313").Value = 150
$ws.Range("Q313").Value = 5
$ws.Range("R313").Value = "Hortaliza"
